# Update Name of Algo
# Applies updated numeric results for the RandomForest imputation output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.019899999999992
$ws.Range("B7").Value = 4.726000000000003
$ws.Range("A10").Value = -21.93579999999999
$ws.Range("A12").Value = -21.54529999999999
$ws.Range("B15").Value = 4.661899999999995
$ws.Range("A18").Value = -22.46480000000001
$ws.Range("D18").Value = -8.276499999999993
$ws.Range("D19").Value = -8.826299999999993
$ws.Range("B20").Value = 9.410199999999998
$ws.Range("D27").Value = -9.095000000000001
$ws.Range("B29").Value = 4.840100000000003
$ws.Range("B30").Value = 5.851400000000002
$ws.Range("B31").Value = 5.509300000000001
$ws.Range("A37").Value = -20.18960000000001
$ws.Range("B40").Value = 9.149499999999994
$ws.Range("D42").Value = -8.867399999999996
$ws.Range("D44").Value = -7.726999999999999
$ws.Range("D47").Value = -7.748200000000002
$ws.Range("A55").Value = -21.8236
$ws.Range("D58").Value = -8.443499999999997
$ws.Range("A68").Value = -21.50500000000001
$ws.Range("B68").Value = 4.688399999999999
$ws.Range("D73").Value = -7.741099999999999
$ws.Range("B76").Value = 5.878799999999999
$ws.Range("A77").Value = -20.57119999999999
$ws.Range("A78").Value = -19.94769999999998
$ws.Range("B87").Value = 4.874599999999996
$ws.Range("B88").Value = 4.5182
$ws.Range("D95").Value = -8.266399999999994
$ws.Range("B96").Value = 5.044900000000005
$ws.Range("B98").Value = 5.566200000000001
$ws.Range("B101").Value = 9.286799999999992
$ws.Range("D101").Value = -7.802299999999998
$ws.Range("B102").Value = 9.020100000000008
